# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

# Update the conversion text block on Hoja1!A1
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 10.26 = 42923.85 pesos"), "1000 Bs = 10.05 = 42995.58 pesos"
$text = $text -replace [regex]::Escape("42923.85 pesos = 10.19 = 946.02 Bs"), "42995.58 pesos = 9.99 = 927.56 Bs"
$cellA1.Value2 = $text

# Update tasas rates on the "tasas" sheet
$wsTasas.Range("N10").Value = 99.5
$wsTasas.Range("O10").Value = 4278.06
$wsTasas.Range("N12").Value = 4303
